$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.685.38'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '2.487.91'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'531.70"
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').Value = "'135.85"
$ws.Range('E6').Value = '  -3.09%  '
$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').Value = '2.507.17'
$ws.Range('E9').Value = '  -1.10%  '
$ws.Range('D10').Value = "'0.100"
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').Value = "'5.34"
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').Value = "'0.346"
$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('D14').Value = '2.933.22'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = "'23.13"
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = '58.632.65'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '2.503.24'
$ws.Range('E18').Value = '  -2.00%  '
$ws.Range('D19').Value = "'11.01"
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('D20').Value = "'4.23"
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').Value = "'323.68"
$ws.Range('E21').Value = '  +0.80%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = "'5.84"
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('D24').Value = "'64.06"
$ws.Range('E24').Value = '  +2.55%  '
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').Value = "'0.996"
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0769'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = "'6.65"
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('E31').Value = '  -1.88%  '
$ws.Range('D32').Value = "'166.89"
$ws.Range('E32').Value = '  +3.56%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = "'0.998"
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = "'1.14"
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('E35').Value = '  -5.22%  '
$ws.Range('D36').Value = "'18.43"
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').Value = "'4.06"
$ws.Range('E37').Value = '  -4.05%  '
$ws.Range('D38').Value = "'1.55"
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('D39').Value = "'36.71"
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('D40').Value = "'0.807"
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('E41').Value = '  -1.12%  '
$ws.Range('D42').Value = "'5.24"
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('D43').Value = "'276.88"
$ws.Range('E43').Value = '  -3.50%  '
$ws.Range('D44').Value = "'0.996"
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = "'10.86"
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').Value = "'128.03"
$ws.Range('E47').Value = '  +3.12%  '
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('D49').Value = "'0.0511"
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').Value = "'17.22"
$ws.Range('E51').Value = '  -1.75%  '
